$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as TEXT even if it looks like a number,
# by forcing the cell's number format to Text ("@") first.
# (used only for D-column values that would otherwise be auto-parsed
# as numbers by Excel, to match the source inlineStr text cells)

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.098.52"
$ws.Range("E2").Value = "  -3.14%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.516.40"
$ws.Range("E3").Value = "  -4.94%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.93"
$ws.Range("E5").Value = "  -0.77%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.07"
$ws.Range("E6").Value = "  -3.83%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.511.16"
$ws.Range("E7").Value = "  -4.88%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.31%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.09%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -5.65%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.78"
$ws.Range("E11").Value = "  -0.84%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -4.53%  "

# Row 13 - Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.00"
$ws.Range("E13").Value = "  -4.21%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -4.59%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.084.19"
$ws.Range("E15").Value = "  -4.87%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -5.40%  "

# Row 17 - BitcoinCash
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "624.41"

# Row 18 - was WrappedBTC, now WrappedEther
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.525.17"
$ws.Range("E18").Value = "  -4.80%  "

# Row 19 - was WrappedEther, now WrappedBTC
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "69.064.94"
$ws.Range("E19").Value = "  -3.28%  "

# Row 21 - Polygon (unchanged label/url)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.40"
$ws.Range("E21").Value = "  -3.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.16"
$ws.Range("E22").Value = "  -3.78%  "

# Row 23
$ws.Range("E23").Value = "  -6.18%  "

# Row 24
$ws.Range("E24").Value = "  -8.93%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.41"
$ws.Range("E25").Value = "  -4.56%  "

# Row 26
$ws.Range("E26").Value = "  -4.58%  "

# Row 27
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
$ws.Range("E28").Value = "  -6.86%  "

# Row 29
$ws.Range("E29").Value = "  -9.74%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.61"
$ws.Range("E30").Value = "  -7.38%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.14"
$ws.Range("E31").Value = "  -8.22%  "

# Row 32
$ws.Range("E32").Value = "  -7.17%  "

# Row 33
$ws.Range("E33").Value = "  -7.43%  "

# Row 34
$ws.Range("E34").Value = "  -6.74%  "

# Row 35 - Bittensor
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "634.23"
$ws.Range("E35").Value = "  +8.25%  "

# Row 36
$ws.Range("E36").Value = "  -4.04%  "

# Row 37
$ws.Range("E37").Value = "  -5.74%  "

# Row 38 - was OKB, now dogwifhat
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.41"
$ws.Range("E38").Value = "  -16.41%  "

# Row 39 - was dogwifhat, now OKB
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.70"
$ws.Range("E39").Value = "  -3.62%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  +0.12%  "

# Row 41 - VeChain
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0445"
$ws.Range("E41").Value = "  -2.26%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -6.16%  "

# Row 43 - Maker
$ws.Range("D43").Value = "3.376.60"
$ws.Range("E43").Value = "  -8.30%  "

# Row 44 - TheGraph
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.327"
$ws.Range("E44").Value = "  -6.95%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.84"
$ws.Range("E45").Value = "  -7.75%  "

# Row 46 - PEPE
$ws.Range("D46").Value = "0.0₃0688"
$ws.Range("E46").Value = "  -10.20%  "

# Row 47 - Fetch.AI
$ws.Range("E47").Value = "  -7.52%  "

# Row 48 - ThetaToken
$ws.Range("E48").Value = "  -4.92%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -2.72%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +14.63%  "

# Row 51 - Monero
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.07"
$ws.Range("E51").Value = "  -2.61%  "
